$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model names for rows 2..26 (row 1 is the header)
$names = @(
    "model_11_7_0",
    "model_11_7_22",
    "model_11_7_21",
    "model_11_7_20",
    "model_11_7_19",
    "model_11_7_18",
    "model_11_7_17",
    "model_11_7_16",
    "model_11_7_15",
    "model_11_7_14",
    "model_11_7_13",
    "model_11_7_23",
    "model_11_7_12",
    "model_11_7_10",
    "model_11_7_9",
    "model_11_7_8",
    "model_11_7_7",
    "model_11_7_6",
    "model_11_7_5",
    "model_11_7_4",
    "model_11_7_3",
    "model_11_7_2",
    "model_11_7_1",
    "model_11_7_11",
    "model_11_7_24"
)

# Common metric values (same values across all rows) for columns B..I
$values = @(
    0.3494677884409869,
    0.4729774674609699,
    0.3599361836854361,
    0.460666355196278,
    0.7199474573135376,
    0.9884014129638672,
    0.5468775033950806,
    0.780625581741333
)

# Use an already-styled cell (A2, style carried over from the original
# sheet) as the formatting source for any newly created rows, so that we
# reuse the existing style index instead of Excel minting a near-duplicate
# cellXf for the new cells.
$styleSource = $ws.Cells.Item(2, 1)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $cellA = $ws.Cells.Item($row, 1)

    if ($row -gt 14) {
        $styleSource.Copy($cellA)
    }

    $cellA.Value = $names[$i]

    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
